# Apply the "mband noise further attempts" update:
#  - Rename algorithm "mband_normal" (row 5) to "mband_noisetry_01"
#  - Add four new algorithm rows (noisetry_03/05/07/09) with new file paths
#  - Update the numeric metrics for row 5 (and the four new rows) to the
#    new recomputed values
#  - Widen column A slightly to fit the new longer label

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width change (14.6328125 -> 17.54296875, to fit the longer
# algorithm labels such as "mband_noisetry_01")
$ws.Columns.Item(1).ColumnWidth = 16.666666666666668

# --- Row 5: rename algorithm + new file path + refreshed metrics ---
$ws.Range("A5").Value = "mband_noisetry_01"
$ws.Range("B5").Value = "C:\Users\gabi\Documents\University\Uni2025\Investigation\PROJECT-25P85\results\EXP2\spectral\NOISE_ESTIMATION\mband_neural_vad_neural_guided_BANDS4_SPACINGLINEAR_FRAME8ms_NEURAL_W0.1.wav"

$metrics = @(
    7.704754856009604,
    3.3462904559752613,
    66.024497435820607,
    1.1556768900720593,
    4.8330441143863947,
    6.4641561894498256,
    3.2655570396750266,
    1.6875575642659519,
    1.5230599701597196,
    1.3738670664190973,
    1.861685810154585,
    2.4671466761269221,
    1.118705181794186,
    1.3143859067555173,
    1.9912552546632791,
    2.0109212801048337,
    1.5430183344654736,
    0.78524730472187798,
    0.80486734112048897,
    0.48618032404703787,
    0.051653188045543641,
    0.23571173354711331
)

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value = $metrics[$i]
}

# --- Rows 6-9: new algorithm variants, same refreshed metrics ---
$newRows = @(
    @{ Row = 6; Algo = "mband_noisetry_03"; Weight = "0.3" },
    @{ Row = 7; Algo = "mband_noisetry_05"; Weight = "0.5" },
    @{ Row = 8; Algo = "mband_noisetry_07"; Weight = "0.7" },
    @{ Row = 9; Algo = "mband_noisetry_09"; Weight = "0.9" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Range("A$r").Value = $entry.Algo
    $ws.Range("B$r").Value = "C:\Users\gabi\Documents\University\Uni2025\Investigation\PROJECT-25P85\results\EXP2\spectral\NOISE_ESTIMATION\mband_neural_vad_neural_guided_BANDS4_SPACINGLINEAR_FRAME8ms_NEURAL_W$($entry.Weight).wav"
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + "$r").Value = $metrics[$i]
    }
}
